$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks so we can rebuild them cleanly in the correct
# order/positions (row insertion does not re-anchor hyperlink refs in this
# engine, so we repopulate every data row directly instead of Insert()).
$ws.Cells.Hyperlinks.Delete()

# row 2
$ws.Cells.Item(2, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(2, 2).Value2 = "AIを活用した社内備品管理アプリ開発の相談"
$ws.Cells.Item(2, 3).Value2 = "システム開発"
$ws.Cells.Item(2, 4).Value2 = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2, 5).Value2 = "期限情報なし"
$ws.Cells.Item(2, 6).Value2 = "https://www.lancers.jp/work/detail/5465005"
$ws.Cells.Item(2, 7).Value2 = 388
$ws.Cells.Item(2, 8).Value2 = "🔥AI,Ai ◆開発 ◇アプリ"

# row 3
$ws.Cells.Item(3, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(3, 2).Value2 = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Cells.Item(3, 3).Value2 = "システム開発"
$ws.Cells.Item(3, 4).Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(3, 5).Value2 = "期限情報なし"
$ws.Cells.Item(3, 6).Value2 = "https://www.lancers.jp/work/detail/5434128"
$ws.Cells.Item(3, 7).Value2 = 368
$ws.Cells.Item(3, 8).Value2 = "🔥AI,Ai ◆開発"

# row 4
$ws.Cells.Item(4, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(4, 2).Value2 = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Cells.Item(4, 3).Value2 = "システム開発"
$ws.Cells.Item(4, 4).Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4, 5).Value2 = "期限情報なし"
$ws.Cells.Item(4, 6).Value2 = "https://www.lancers.jp/work/detail/5427956"
$ws.Cells.Item(4, 7).Value2 = 310
$ws.Cells.Item(4, 8).Value2 = "🔥AI,Ai"

# row 5
$ws.Cells.Item(5, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(5, 2).Value2 = "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)"
$ws.Cells.Item(5, 3).Value2 = "システム開発"
$ws.Cells.Item(5, 4).Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(5, 5).Value2 = "期限情報なし"
$ws.Cells.Item(5, 6).Value2 = "https://www.lancers.jp/work/detail/5439158"
$ws.Cells.Item(5, 7).Value2 = 303
$ws.Cells.Item(5, 8).Value2 = "🔥AI,Ai"

# row 6
$ws.Cells.Item(6, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(6, 2).Value2 = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Cells.Item(6, 3).Value2 = "システム開発"
$ws.Cells.Item(6, 4).Value2 = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(6, 5).Value2 = "期限情報なし"
$ws.Cells.Item(6, 6).Value2 = "https://www.lancers.jp/work/detail/5405023"
$ws.Cells.Item(6, 7).Value2 = 178
$ws.Cells.Item(6, 8).Value2 = "★bot ◆ツール"

# row 7
$ws.Cells.Item(7, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(7, 2).Value2 = "【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発"
$ws.Cells.Item(7, 3).Value2 = "システム開発"
$ws.Cells.Item(7, 4).Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(7, 5).Value2 = "期限情報なし"
$ws.Cells.Item(7, 6).Value2 = "https://www.lancers.jp/work/detail/5454210"
$ws.Cells.Item(7, 7).Value2 = 108
$ws.Cells.Item(7, 8).Value2 = "◆開発 ◇アプリ"

# row 8
$ws.Cells.Item(8, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(8, 2).Value2 = "GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集"
$ws.Cells.Item(8, 3).Value2 = "システム開発"
$ws.Cells.Item(8, 4).Value2 = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(8, 5).Value2 = "期限情報なし"
$ws.Cells.Item(8, 6).Value2 = "https://www.lancers.jp/work/detail/5457458"
$ws.Cells.Item(8, 7).Value2 = 75
$ws.Cells.Item(8, 8).Value2 = "◆開発"

# row 9
$ws.Cells.Item(9, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(9, 2).Value2 = "Octoparseを使ったスクレイピングシステムの構築"
$ws.Cells.Item(9, 3).Value2 = "システム開発"
$ws.Cells.Item(9, 4).Value2 = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(9, 5).Value2 = "期限情報なし"
$ws.Cells.Item(9, 6).Value2 = "https://www.lancers.jp/work/detail/5465301"
$ws.Cells.Item(9, 7).Value2 = 58
$ws.Cells.Item(9, 8).Value2 = "◆スクレイピング"

# row 10
$ws.Cells.Item(10, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(10, 2).Value2 = "【急募】Appsheetで見積もりアプリを作成してくれる方"
$ws.Cells.Item(10, 3).Value2 = "システム開発"
$ws.Cells.Item(10, 4).Value2 = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(10, 5).Value2 = "期限情報なし"
$ws.Cells.Item(10, 6).Value2 = "https://www.lancers.jp/work/detail/5465442"
$ws.Cells.Item(10, 7).Value2 = 30
$ws.Cells.Item(10, 8).Value2 = "◇アプリ"

# row 11
$ws.Cells.Item(11, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(11, 2).Value2 = "初回 bubbleで構築したサイトの修正対応"
$ws.Cells.Item(11, 3).Value2 = "システム開発"
$ws.Cells.Item(11, 4).Value2 = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(11, 5).Value2 = "期限情報なし"
$ws.Cells.Item(11, 6).Value2 = "https://www.lancers.jp/work/detail/5465187"
$ws.Cells.Item(11, 7).Value2 = 30
$ws.Cells.Item(11, 8).Value2 = "◇サイト"

# row 12
$ws.Cells.Item(12, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(12, 2).Value2 = "【急募】メール問い合わせ時の自動SMS送信システム構築"
$ws.Cells.Item(12, 3).Value2 = "システム開発"
$ws.Cells.Item(12, 4).Value2 = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(12, 5).Value2 = "期限情報なし"
$ws.Cells.Item(12, 6).Value2 = "https://www.lancers.jp/work/detail/5464796"
$ws.Cells.Item(12, 7).Value2 = 33
$ws.Cells.Item(12, 8).ClearContents()

# row 13
$ws.Cells.Item(13, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(13, 2).Value2 = "金融機関の入出金伝票印刷システム構築依頼"
$ws.Cells.Item(13, 3).Value2 = "システム開発"
$ws.Cells.Item(13, 4).Value2 = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(13, 5).Value2 = "期限情報なし"
$ws.Cells.Item(13, 6).Value2 = "https://www.lancers.jp/work/detail/5464833"
$ws.Cells.Item(13, 7).Value2 = 28
$ws.Cells.Item(13, 8).ClearContents()

# row 14
$ws.Cells.Item(14, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(14, 2).Value2 = "【急募】クラウドウェア内製化推進のための技術サポート依頼"
$ws.Cells.Item(14, 3).Value2 = "システム開発"
$ws.Cells.Item(14, 4).Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(14, 5).Value2 = "期限情報なし"
$ws.Cells.Item(14, 6).Value2 = "https://www.lancers.jp/work/detail/5465210"
$ws.Cells.Item(14, 7).Value2 = 25
$ws.Cells.Item(14, 8).ClearContents()

# row 15
$ws.Cells.Item(15, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(15, 2).Value2 = "【退職サポート】退職申請フロー・退会申請フローのチャットボット作成依頼"
$ws.Cells.Item(15, 3).Value2 = "システム開発"
$ws.Cells.Item(15, 4).Value2 = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(15, 5).Value2 = "期限情報なし"
$ws.Cells.Item(15, 6).Value2 = "https://www.lancers.jp/work/detail/5465526"
$ws.Cells.Item(15, 7).Value2 = 18
$ws.Cells.Item(15, 8).ClearContents()

# row 16
$ws.Cells.Item(16, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(16, 2).Value2 = "限定公開 限定公開の仕事"
$ws.Cells.Item(16, 3).Value2 = "システム開発"
$ws.Cells.Item(16, 4).Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(16, 5).Value2 = "期限情報なし"
$ws.Cells.Item(16, 6).Value2 = "https://www.lancers.jp/work/detail/5465372"
$ws.Cells.Item(16, 7).Value2 = 18
$ws.Cells.Item(16, 8).ClearContents()

# row 17
$ws.Cells.Item(17, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(17, 2).Value2 = "RobloxアクションRPG制作(MVP/完成版前提プロジェクト)"
$ws.Cells.Item(17, 3).Value2 = "システム開発"
$ws.Cells.Item(17, 4).Value2 = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(17, 5).Value2 = "期限情報なし"
$ws.Cells.Item(17, 6).Value2 = "https://www.lancers.jp/work/detail/5465063"
$ws.Cells.Item(17, 7).Value2 = 18
$ws.Cells.Item(17, 8).ClearContents()

# row 18
$ws.Cells.Item(18, 1).Value2 = "2026-01-05 18:29:24"
$ws.Cells.Item(18, 2).Value2 = "【準委任】音声データ収集プロジェクトのPM・ディレクター募集"
$ws.Cells.Item(18, 3).Value2 = "システム開発"
$ws.Cells.Item(18, 4).Value2 = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(18, 5).Value2 = "期限情報なし"
$ws.Cells.Item(18, 6).Value2 = "https://www.lancers.jp/work/detail/5465028"
$ws.Cells.Item(18, 7).Value2 = 18
$ws.Cells.Item(18, 8).ClearContents()

# Rebuild hyperlinks on column F for rows 2..18, in order.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5465005") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5434128") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5427956") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5439158") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5405023") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5454210") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5457458") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5465301") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5465442") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5465187") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5464796") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5464833") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5465210") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5465526") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5465372") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5465063") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5465028") | Out-Null

# Column B width: raw xlsx width 41 -> 51 (ColumnWidth property uses a
# different character-width unit with a constant ~0.83 padding offset).
$ws.Columns.Item(2).ColumnWidth = 50.17

"done"